$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the hyphen with a middle dot (·) in the existing textual sprint-item
# identifiers so Excel no longer has a reason to treat them as dates.
$ws.Range("A4").Value = "0·01"
$ws.Range("A8").Value = "0·04"

# A10 and A13 used to hold the literal dates 43466 / 43467 (1-01 / 2-01 got
# auto-converted to dates by Excel). Clear the date number format back to
# General before writing the correct textual identifiers.
$ws.Range("A10").NumberFormat = "general"
$ws.Range("A10").Value = "1·01"

$ws.Range("A13").NumberFormat = "general"
$ws.Range("A13").Value = "2·01"
